# The workbook originally had its data starting at A1 on every sheet.
# Restore the "original" layout where every sheet's table/data block is
# shifted down by one row and right by one column (now starting at B2),
# by inserting a blank row above row 1 and a blank column to the left of
# column A on each worksheet, then resizing the Excel Tables (ListObjects)
# to their new location.

$wb = $excel.ActiveWorkbook

# --- Sheet "user" (Table1: username/password/role) ---
$wsUser = $wb.Worksheets.Item("user")
$wsUser.Rows.Item(1).Insert()
$wsUser.Columns.Item(1).Insert()
$wsUser.ListObjects.Item(1).Resize($wsUser.Range("B2:D5"))
$wsUser.Range("B2:D5").Select()

# --- Sheet "product" (Table2: productId/name/description/imageLocation/price) ---
$wsProduct = $wb.Worksheets.Item("product")
$wsProduct.Rows.Item(1).Insert()
$wsProduct.Columns.Item(1).Insert()
$wsProduct.ListObjects.Item(1).Resize($wsProduct.Range("B2:F8"))
$wsProduct.Range("F22").Select()

# --- Sheet "saleTransaction" (Table3: transactionId/ProductId/amount) ---
$wsSale = $wb.Worksheets.Item("saleTransaction")
$wsSale.Rows.Item(1).Insert()
$wsSale.Columns.Item(1).Insert()
$wsSale.ListObjects.Item(1).Resize($wsSale.Range("B2:D8"))
$wsSale.Range("B8").Select()

# --- Sheet "SaleOrder" (plain range, no table) ---
$wsOrder = $wb.Worksheets.Item("SaleOrder")
$wsOrder.Rows.Item(1).Insert()
$wsOrder.Columns.Item(1).Insert()

# Keep "SaleOrder" the active tab, matching the workbook's saved state.
$wsOrder.Activate()
$wsOrder.Range("C4").Select()
